$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new header row at the top; existing rows shift down by one.
$ws.Rows.Item(1).Insert()

$ws.Range("A1").Value = "Município"
$ws.Range("B1").Value = "Casos"
$ws.Range("C1").Value = "Óbitos"

# Fix spelling correction that landed on the shifted row (was "mogi mirim").
$ws.Range("A125").Value = "moji mirim"

$headerRange = $ws.Range("A1:C1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
